# Automatic update of files.
# Column C ("Förändrad") on rows 2-16 gets bumped from date serial 45243
# (2023-11-13) to 45244 (2023-11-14) for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
